$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue 'D2' '27.881.50'
$ws.Range('E2').Value = '  +2.74%  '
Set-TextValue 'D3' '1.663.52'
$ws.Range('E3').Value = '  -0.88%  '
$ws.Range('E4').Value = '  +0.34%  '
Set-TextValue 'D5' '215.05'
$ws.Range('E5').Value = '  +0.39%  '
Set-TextValue 'D6' '0.516'
$ws.Range('E6').Value = '  -0.32%  '
$ws.Range('E7').Value = '  +0.26%  '
Set-TextValue 'D8' '23.51'
$ws.Range('E9').Value = '  +0.69%  '
Set-TextValue 'D10' '0.0620'
$ws.Range('E10').Value = '  -0.32%  '
Set-TextValue 'D11' '0.0883'
$ws.Range('E11').Value = '  -0.80%  '
Set-TextValue 'D12' '1.905.68'
$ws.Range('E12').Value = '  -0.53%  '
Set-TextValue 'D13' '1.680.92'
$ws.Range('E13').Value = '  +0.16%  '
Set-TextValue 'D14' '4.14'
$ws.Range('E14').Value = '  -1.40%  '
Set-TextValue 'D15' '0.550'
$ws.Range('E15').Value = '  -0.08%  '
Set-TextValue 'D16' '66.22'
$ws.Range('E16').Value = '  -0.55%  '
Set-TextValue 'D17' '247.73'
$ws.Range('E17').Value = '  +5.12%  '
Set-TextValue 'D18' '27.854.21'
$ws.Range('E18').Value = '  +2.77%  '
Set-TextValue 'D19' '0.0₃0734'
$ws.Range('E19').Value = '  -0.93%  '
Set-TextValue 'D20' '7.57'
$ws.Range('E20').Value = '  -4.16%  '
Set-TextValue 'D21' '1.01'
$ws.Range('E21').Value = '  +0.43%  '
Set-TextValue 'D22' '4.47'
$ws.Range('E22').Value = '  -1.54%  '
Set-TextValue 'D23' '9.16'
$ws.Range('E23').Value = '  -3.72%  '
Set-TextValue 'D24' '2.04'
$ws.Range('E24').Value = '  -2.34%  '
Set-TextValue 'D25' '146.78'
$ws.Range('E25').Value = '  -0.29%  '
Set-TextValue 'D26' '7.24'
$ws.Range('E26').Value = '  -2.48%  '
Set-TextValue 'D27' '16.23'
$ws.Range('E27').Value = '  -0.68%  '
$ws.Range('B28').Value = 'BinanceUSD'
$ws.Range('C28').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue 'D28' '1.01'
$ws.Range('E28').Value = '  +0.38%  '
$ws.Range('B29').Value = 'Stellar'
$ws.Range('C29').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue 'D29' '0.112'
$ws.Range('E29').Value = '  -0.69%  '
Set-TextValue 'D30' '1.24'
$ws.Range('E30').Value = '  +6.56%  '
Set-TextValue 'D31' '0.0500'
$ws.Range('E31').Value = '  -0.35%  '
Set-TextValue 'D32' '3.34'
$ws.Range('E32').Value = '  -0.71%  '
Set-TextValue 'D33' '3.13'
$ws.Range('E33').Value = '  -3.39%  '
Set-TextValue 'D34' '1.415.79'
$ws.Range('E34').Value = '  -8.18%  '
Set-TextValue 'D35' '1.57'
$ws.Range('E35').Value = '  -5.05%  '
Set-TextValue 'D36' '2.40'
$ws.Range('E36').Value = '  +0.54%  '
$ws.Range('E37').Value = '  -1.25%  '
Set-TextValue 'D38' '0.579'
$ws.Range('E38').Value = '  -4.57%  '
$ws.Range('E39').Value = '  -1.17%  '
$ws.Range('E40').Value = '  -2.66%  '
$ws.Range('B41').Value = 'PaxDollar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue 'D41' '1.01'
$ws.Range('E41').Value = '  +0.38%  '
$ws.Range('B42').Value = 'Aave'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 'D42' '68.98'
$ws.Range('E42').Value = '  -0.87%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue 'D43' '5.45'
$ws.Range('E43').Value = '  -6.02%  '
$ws.Range('B44').Value = 'MXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue 'D44' '2.23'
$ws.Range('E44').Value = '  -0.95%  '
$ws.Range('B45').Value = 'TrustWalletToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue 'D45' '0.794'
$ws.Range('E45').Value = '  +1.99%  '
$ws.Range('B46').Value = 'RocketPoolETH'
$ws.Range('C46').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
Set-TextValue 'D46' '1.809.94'
$ws.Range('E46').Value = '  -0.66%  '
$ws.Range('E47').Value = '  +4.69%  '
Set-TextValue 'D48' '88.38'
$ws.Range('E48').Value = '  -1.67%  '
Set-TextValue 'D49' '0.0₆0107'
$ws.Range('E49').Value = '  -3.98%  '
Set-TextValue 'D50' '0.101'
$ws.Range('E50').Value = '  -2.49%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue 'D51' '0.0509'
$ws.Range('E51').Value = '  -0.26%  '
